$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inputs
$ws.Range("B3").Value = 4
$ws.Range("B4").Value = 18
$ws.Range("B7").Value = 10
$ws.Range("B8").Value = 5

# Formulas
$ws.Range("B5").Formula = "=B4-B3+1"
$ws.Range("E3").Formula = "=1/B5"
$ws.Range("E5").Formula = "=1/B5+1/B5"
$ws.Range("F5").Formula = "=2*E3"
$ws.Range("E7").Formula = "=E3*E3"
$ws.Range("B10").Formula = "=RAND()"

$ws.Range("B11").Select()
